$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.286.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.890.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.742"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000354"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.12"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.50"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.505.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.884.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.281.91"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.03%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.46"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +16.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.65"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.69"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "715.25"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.90"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0881"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.07"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.03"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.400"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +15.70%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.07%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.89%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0497"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.69%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.90%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.11%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +27.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.11"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.04%  "
